# Apply the scheduled cryptos.xlsx price/volume refresh
# Commit: "Updated cryptos list on Wed Feb 14 23:45:08 UTC 2024 with GitHub Actions"
#
# D/E columns hold plain text (not real numbers), e.g. "117.10" must keep its
# trailing zero and "0.130" must not become "0.13". Excel's COM Range.Value
# setter auto-coerces any numeric-looking string to a number, so for D-column
# values that would parse as a number we force the Text number format first
# (mirrors typing into a cell pre-formatted as Text) and then write the literal
# string. Values that already resist numeric parsing (e.g. thousands-grouped
# "51.826.49", the subscript-containing "0.0₃0975") and the percent-change
# column (always has stray spaces + a % sign) are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.826.49'
$ws.Range("E2").Value = '  +4.55%  '
$ws.Range("D3").Value = '2.780.42'
$ws.Range("E3").Value = '  +5.53%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '117.10'
$ws.Range("E5").Value = '  +4.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '333.95'
$ws.Range("E6").Value = '  +2.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.538'
$ws.Range("E7").Value = '  +2.39%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +5.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.29'
$ws.Range("E10").Value = '  +6.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0855'
$ws.Range("E11").Value = '  +5.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.27'
$ws.Range("E12").Value = '  +2.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.130'
$ws.Range("E13").Value = '  +2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.65'
$ws.Range("E14").Value = '  +4.16%  '
$ws.Range("D15").Value = '3.212.17'
$ws.Range("E15").Value = '  +5.33%  '
$ws.Range("D16").Value = '2.779.12'
$ws.Range("E16").Value = '  +5.40%  '
$ws.Range("E17").Value = '  +4.89%  '
$ws.Range("D18").Value = '51.812.00'
$ws.Range("E18").Value = '  +4.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.29'
$ws.Range("E19").Value = '  +11.09%  '
$ws.Range("E20").Value = '  +5.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.88'
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("D22").Value = '0.0₃0975'
$ws.Range("E22").Value = '  +3.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '280.44'
$ws.Range("E23").Value = '  +3.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.94'
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("E25").Value = '  +6.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.90'
$ws.Range("E26").Value = '  +2.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.24'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("E29").Value = '  +0.96%  '
$ws.Range("E30").Value = '  +3.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.32'
$ws.Range("E31").Value = '  +1.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.20'
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0824'
$ws.Range("E34").Value = '  +1.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.16'
$ws.Range("E35").Value = '  +1.00%  '
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.11'
$ws.Range("E37").Value = '  +3.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.03'
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.29'
$ws.Range("E39").Value = '  +6.05%  '
$ws.Range("E40").Value = '  +10.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '127.71'
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.55'
$ws.Range("E42").Value = '  +20.34%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.32'
$ws.Range("E43").Value = '  +6.12%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.32'
$ws.Range("E44").Value = '  +7.91%  '
$ws.Range("E45").Value = '  +2.81%  '
$ws.Range("D46").Value = '2.090.07'
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  +4.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.54'
$ws.Range("E49").Value = '  +6.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '60.85'
$ws.Range("E50").Value = '  +3.05%  '
$ws.Range("E51").Value = '  -0.44%  '
